$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.640.94'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '2.546.30'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.15'
$ws.Range('E5').Value = '  -1.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.37'
$ws.Range('E6').Value = '  +5.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.569'
$ws.Range('E7').Value = '  -1.25%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.08'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  -1.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.40'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').Value = '2.936.61'
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.97'
$ws.Range('E15').Value = '  +6.01%  '
$ws.Range('D16').Value = '2.537.95'
$ws.Range('E16').Value = '  -0.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.834'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('D18').Value = '42.664.22'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.82'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.33'
$ws.Range('E20').Value = '  -2.10%  '
$ws.Range('D21').Value = '0.0₃0953'
$ws.Range('E21').Value = '  -1.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.10'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '244.17'
$ws.Range('E23').Value = '  -3.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.92'
$ws.Range('E24').Value = '  -1.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.06'
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.47'
$ws.Range('E26').Value = '  -1.51%  '
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('E29').Value = '  -1.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.10'
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.81'
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.71'
$ws.Range('E32').Value = '  -3.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.75'
$ws.Range('E33').Value = '  +11.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0804'
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.63'
$ws.Range('E35').Value = '  -3.13%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.06'
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.24'
$ws.Range('E37').Value = '  -3.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.42'
$ws.Range('E38').Value = '  -4.67%  '
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('E40').Value = '  -0.81%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.33'
$ws.Range('E41').Value = '  +2.49%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.20'
$ws.Range('E42').Value = '  +9.87%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.32'
$ws.Range('E44').Value = '  +1.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0298'
$ws.Range('E45').Value = '  -2.00%  '
$ws.Range('D46').Value = '1.973.29'
$ws.Range('E46').Value = '  -1.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.91'
$ws.Range('E47').Value = '  -2.03%  '
$ws.Range('D48').Value = '2.794.54'
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '81.16'
$ws.Range('E49').Value = '  -4.13%  '
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.25'
$ws.Range('E51').Value = '  -2.03%  '

# Strip the temporary text-number-format override back off so the
# cells end up with the workbook's default (unstyled) appearance,
# same as every other data cell in the sheet.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
